$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spreadsheet headers: rename "Sample ID" -> "Specimen_Number" and "T/N" -> "SAMPLE_TYPE"
# A1's original cell style uses a quote-prefix (text forced via leading apostrophe), so
# preserve that by assigning via Formula with a leading apostrophe instead of Value.
$ws.Range("A1").Formula = "'Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Update the selected/active cell to F2
$ws.Range("F2").Select()
